$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 325.125
$ws.Range("J17").Value = 325.125
$ws.Range("L17").Value = 975.375
$ws.Range("N17").Value = -1311.375

# Row 111
$ws.Range("H111").Value = 2118.2
$ws.Range("I111").Value = 1842.6666
$ws.Range("J111").Value = 2531.5
$ws.Range("K111").Value = 5527.9998
$ws.Range("L111").Value = 7594.5
$ws.Range("M111").Value = -2460.9998
$ws.Range("N111").Value = -13728.5

# Row 113
$ws.Range("H113").Value = 16499.25
$ws.Range("I113").Value = 16499.25
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 16499.25
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -13245.25
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 141
$ws.Range("I5").Value = 411
$ws.Range("K5").Value = 411
$ws.Range("M5").Value = -299

# Row 61
$ws.Range("H61").Value = 1102.4
$ws.Range("I61").Value = 1124.5
$ws.Range("K61").Value = 1124.5
$ws.Range("M61").Value = -912.5

# Row 74
$ws.Range("H74").Value = 1110.1428
$ws.Range("I74").Value = 1024.3334
$ws.Range("K74").Value = 1024.3334
$ws.Range("M74").Value = -150.3334

# Row 77
$ws.Range("H77").Value = 1110.1428
$ws.Range("I77").Value = 1024.3334
$ws.Range("K77").Value = 5121.666999999999
$ws.Range("M77").Value = -753.6669999999995

# Row 132
$ws.Range("H132").Value = 2598
$ws.Range("I132").Value = 2598
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7794
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5264
$ws.Range("N132").ClearContents()

# Row 136
$ws.Range("H136").Value = 1102.4
$ws.Range("I136").Value = 1124.5
$ws.Range("K136").Value = 3373.5
$ws.Range("M136").Value = -823.5

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 141
$ws.Range("I4").Value = 411
$ws.Range("K4").Value = 411
$ws.Range("M4").Value = -296

# Row 22
$ws.Range("H22").Value = 360.66666
$ws.Range("I22").Value = 360.66666
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 360.66666
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -187.66666
$ws.Range("N22").ClearContents()

# Row 102
$ws.Range("H102").Value = 9778
$ws.Range("I102").Value = 9778
$ws.Range("K102").Value = 9778
$ws.Range("M102").Value = -6533

# Row 134
$ws.Range("H134").Value = 3670.4119
$ws.Range("I134").Value = 3946.8
$ws.Range("K134").Value = 11840.4
$ws.Range("M134").Value = -9305.400000000001

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2141.6843
$ws.Range("I31").Value = 1664.6
$ws.Range("J31").Value = 3930.75
$ws.Range("K31").Value = 1664.6
$ws.Range("L31").Value = 3930.75
$ws.Range("M31").Value = -1369.6
$ws.Range("N31").Value = -4520.75

# Row 34
$ws.Range("H34").Value = 2141.6843
$ws.Range("I34").Value = 1664.6
$ws.Range("J34").Value = 3930.75
$ws.Range("K34").Value = 1664.6
$ws.Range("L34").Value = 3930.75
$ws.Range("M34").Value = -1462.6
$ws.Range("N34").Value = -4334.75

# Row 51
$ws.Range("H51").Value = 20000
$ws.Range("J51").Value = 20000
$ws.Range("L51").Value = 20000
$ws.Range("N51").Value = -21472

# Row 61
$ws.Range("H61").Value = 20000
$ws.Range("J61").Value = 20000
$ws.Range("L61").Value = 20000
$ws.Range("N61").Value = -20696

# Row 95
$ws.Range("H95").Value = 12054
$ws.Range("J95").Value = 12054
$ws.Range("L95").Value = 12054
$ws.Range("N95").Value = -17546

$ws = $wb.Worksheets.Item("CUL")
# Row 81
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()

# Row 84
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()

# Row 86
$ws.Range("H86").Value = 1054
$ws.Range("J86").Value = 949.5
$ws.Range("L86").Value = 2848.5
$ws.Range("N86").Value = -5220.5

# Row 87
$ws.Range("H87").Value = 500
$ws.Range("I87").Value = 500
$ws.Range("K87").Value = 1500
$ws.Range("M87").Value = -252

# Row 89
$ws.Range("H89").Value = 1054
$ws.Range("J89").Value = 949.5
$ws.Range("L89").Value = 8545.5
$ws.Range("N89").Value = -20401.5

# Row 90
$ws.Range("H90").Value = 500
$ws.Range("I90").Value = 500
$ws.Range("K90").Value = 4500
$ws.Range("M90").Value = 1740

# Row 119
$ws.Range("H119").Value = 647.7143
$ws.Range("I119").Value = 647.7143
$ws.Range("K119").Value = 1943.1429
$ws.Range("M119").Value = 2894.8571

# Row 131
$ws.Range("H131").Value = 987.1875

# Row 140
$ws.Range("H140").Value = 909.6
$ws.Range("I140").Value = 909.6
$ws.Range("K140").Value = 2728.8
$ws.Range("M140").Value = 2451.2

$ws = $wb.Worksheets.Item("GSM")
# Row 19
$ws.Range("H19").Value = 1666.3334
$ws.Range("I19").Value = 2332.6667
$ws.Range("K19").Value = 2332.6667
$ws.Range("M19").Value = -2044.6667

# Row 54
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("N54").ClearContents()

# Row 97
$ws.Range("H97").Value = 2678.5715
$ws.Range("I97").Value = 2473
$ws.Range("K97").Value = 2473
$ws.Range("M97").Value = -1977

# Row 102
$ws.Range("H102").Value = 1842.3334
$ws.Range("I102").Value = 1426
$ws.Range("J102").Value = 2675
$ws.Range("K102").Value = 1426
$ws.Range("L102").Value = 2675
$ws.Range("M102").Value = 196
$ws.Range("N102").Value = -5919

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 6098.5713
$ws.Range("I16").Value = 5948
$ws.Range("K16").Value = 5948
$ws.Range("M16").Value = -5778

# Row 93
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 1000
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -3496

# Row 132
$ws.Range("H132").Value = 4434.778
$ws.Range("J132").Value = 6253.5
$ws.Range("L132").Value = 18760.5
$ws.Range("N132").Value = -23820.5

# Row 136
$ws.Range("H136").Value = 8049
$ws.Range("I136").Value = 8054.4443
$ws.Range("K136").Value = 24163.3329
$ws.Range("M136").Value = -21613.3329

$ws = $wb.Worksheets.Item("WVR")
# Row 33
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

# Row 36
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

# Row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

# Row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
